$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update cell values (rows 6-8 content shuffled up) ---
$ws.Range("A6").Value = "debashree.p@insync.co.inxxxxxxx"
$ws.Range("B6").Value = "Efgh.1234"
$ws.Range("A7").Value = "rio1@yopmail.com"
$ws.Range("B7").Value = ""
$ws.Range("A8").Value = ""

# --- Rebuild hyperlinks to match the new layout ---
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A2"), "mailto:debashree.p@insync.co.in")
$ws.Hyperlinks.Add($ws.Range("A4"), "mailto:debashree.p@insync.co.in")
$ws.Hyperlinks.Add($ws.Range("A3"), "mailto:rio1@yopmail.com")
$ws.Hyperlinks.Add($ws.Range("A5"), "mailto:debashree.p@insync.co.in")
$ws.Hyperlinks.Add($ws.Range("A6"), "mailto:debashree.p@insync.co.inxxxxxxx")
$ws.Hyperlinks.Add($ws.Range("A7"), "mailto:rio1@yopmail.com")

# Hyperlinks.Add() reformats the target cells with a generic "Hyperlink" style;
# restore the worksheet's original hyperlink-cell look (style carried by A8,
# which keeps its look untouched by the edit) on every linked cell.
$ws.Range("A8").Copy()
$ws.Range("A2:A7").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Selection left by the editor ---
$ws.Range("A13").Select()
